$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price (D) column so numeric-looking strings
# (e.g. "29.265.09", "1.000") are preserved as text instead of being
# parsed into numbers, matching the original inline-string cell types.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.265.09"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "1.839.88"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "243.40"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "0.6866"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.3028"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").Value = "0.07513"
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("D10").Value = "23.21"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").Value = "0.07691"
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("D12").Value = "1.836.77"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "5.083"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "0.6859"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "88.56"
$ws.Range("E15").Value = "  -4.51%  "
$ws.Range("D16").Value = "6.268"
$ws.Range("E16").Value = "  -3.44%  "
$ws.Range("D17").Value = "29.277.31"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").Value = "0.000008216"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "2.089.39"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").Value = "232.09"
$ws.Range("E20").Value = "  -3.84%  "
$ws.Range("D21").Value = "12.59"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "7.470"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "0.1460"
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("D26").Value = "159.75"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").Value = "8.824"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").Value = "18.12"
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").Value = "1.518"
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("D30").Value = "4.244"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").Value = "4.150"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").Value = "1.203"
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("D33").Value = "0.05165"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").Value = "0.7710"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("D35").Value = "1.835"
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("D36").Value = "1.138"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").Value = "2.674"
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("D38").Value = "1.300.50"
$ws.Range("E38").Value = "  +2.06%  "
$ws.Range("D39").Value = "0.01842"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").Value = "2.698"
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").Value = "0.9448"
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("D42").Value = "105.05"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("D43").Value = "5.783"
$ws.Range("E43").Value = "  -5.85%  "
$ws.Range("D44").Value = "0.9998"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "9.662"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "1.986.55"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("D47").Value = "0.5197"
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("D48").Value = "64.54"
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("D49").Value = "1.769"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("D50").Value = "0.00000000119"
$ws.Range("E50").Value = "  -2.94%  "
$ws.Range("D51").Value = "0.05922"
$ws.Range("E51").Value = "  +0.85%  "

# Row 42/43: coin name + link swap (FraxShare <-> Quant)
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"

# Reset the Price column style back to its original (default) style so no
# extra style index gets attached to the cells themselves.
$ws.Range("D2:D51").Style = "Normal"
